$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C column) date value from 45184 to 45186 for every
#    data row (rows 2 through 452).
$ws.Range("C2:C452").Value2 = 45186

# 2) Add the "friendly name" (the value from column A, e.g. "A 69213-2018") as
#    the second argument of every HYPERLINK(...) formula found in columns
#    S, T, U, V, W, X and Y, for the data rows that contain such formulas
#    (rows 2 through 26).
$cols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le 26; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2

    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f.StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
                $inner = $f.Substring(1, $f.Length - 2)
                $cell.Formula = "=" + $inner + ', "' + $label + '")'
            }
        }
    }
}
